# feat(infra): agregar Config Server al monorepo de microservicios
#
# Highlight (in red, RGB EE0000) the two bullet items under the
# "Base del Proyecto" section that describe the Git repo layout and the
# new Config Server, calling out the infra change described in the
# commit message:
#   - "Crear un repositorio Git (puede ser monorepo con carpetas ...)"
#   - "Configurar un Config Server con un repositorio de configuración."
#
# wdColor is packed as 0x00BBGGRR, so pure red (OOXML w:val="EE0000",
# i.e. R=0xEE, G=0x00, B=0x00) is just the integer 0xEE (238).

$d = $word.ActiveDocument
$wdRed = 238

$gitRepoPara = $null
$configServerPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text

    if ($null -eq $gitRepoPara -and $text -like "Crear un repositorio Git*") {
        $gitRepoPara = $para
    }
    elseif ($null -eq $configServerPara -and $text -like "Configurar un Config Server*") {
        $configServerPara = $para
    }
}

if ($gitRepoPara) {
    $gitRepoPara.Range.Font.Color = $wdRed
}

if ($configServerPara) {
    $configServerPara.Range.Font.Color = $wdRed
}
